# Regenerate merged AHB files
# Renames header columns from *_old/*_new to *_FV2410/*_FV2504,
# adds an Excel Table over the data range, and freezes the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2410 = "_FV2410"
$fv2504 = "_FV2504"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count()
$lastCol = $usedRange.Columns.Count()

# Rename header row values (row 1): "*_old" -> "*_FV2410", "*_new" -> "*_FV2504"
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith($oldSuffix)) {
            $base = $val.Substring(0, $val.Length - $oldSuffix.Length)
            $cell.Value = $base + $fv2410
        } elseif ($val.EndsWith($newSuffix)) {
            $base = $val.Substring(0, $val.Length - $newSuffix.Length)
            $cell.Value = $base + $fv2504
        }
    }
}

# Determine the full data range (A1:U60)
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

# Add an Excel Table (ListObject) over the data range
$listObj = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$listObj.Name = "Table1"
# Match source workbook: no explicit named table style applied
$listObj.TableStyle = ""

# Freeze the header row (freeze panes at A2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
